$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows of Mac-Address (user-machine) data, following the existing pattern
$newRows = @(
    @(10001, 110030, 10030),
    @(10001, 110031, 10031)
)

$startRow = 31
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
}

# Update the view/selection to match the final state: scrolled to row 25, active cell C29
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("C29").Select()
